$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
[double]$x = 1.1
$ws.Range("A1").Value = $x
$v = $ws.Range("A1").Value()
Write-Host "A1=$v"

$ws.Range("A2").Value = [double]1.1
$v2 = $ws.Range("A2").Value()
Write-Host "A2=$v2"
